$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5858
$ws1.Range("F6").Value = 5178
$ws1.Range("F11").Value = 218

# Sheet "全部类型" (all types) - same underlying rows, mirror the update
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5858
$ws4.Range("F6").Value = 5178
$ws4.Range("F11").Value = 218
